$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ridge)
$ws.Range("B2").Value = 0.985979005791894
$ws.Range("C2").Value = 0.9237427832909219
$ws.Range("D2").Value = 0.06223622250097216
$ws.Range("E2").Value = 3.883903536348642
$ws.Range("F2").Value = 9.057732799156215
$ws.Range("G2").Value = "{'solver': 'saga', 'alpha': 1.0}"
$ws.Range("H2").Value = 3.96

# Row 3 (Lasso)
$ws.Range("B3").Value = 0.986294271471534
$ws.Range("C3").Value = 0.9247897503465051
$ws.Range("D3").Value = 0.06150452112502891
$ws.Range("E3").Value = 3.839989993488935
$ws.Range("F3").Value = 8.99533921105319
$ws.Range("G3").Value = "{'selection': 'random', 'alpha': 0.01}"
$ws.Range("H3").Value = 0.79

# Row 4 (ElasticNet)
$ws.Range("H4").Value = 2.07

# Row 5 (SVR)
$ws.Range("H5").Value = 5.13

# Row 6 (KNN Regressor)
$ws.Range("H6").Value = 6.55

# Row 7 (Decision Tree)
$ws.Range("H7").Value = 2.15

# Row 8 (PLSRegression)
$ws.Range("H8").Value = 0.5600000000000001
